$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 25148.75
$ws.Range("I28").Value = 25148.75
$ws.Range("K28").Value = 25148.75
$ws.Range("M28").Value = -24663.75

$ws.Range("H32").Value = 2261.5386
$ws.Range("J32").Value = 2666.6667
$ws.Range("L32").Value = 2666.6667
$ws.Range("N32").Value = -3318.6667

$ws.Range("H69").Value = 3546
$ws.Range("I69").Value = 3244.3333
$ws.Range("K69").Value = 9732.999899999999
$ws.Range("M69").Value = -8858.999899999999

$ws.Range("H72").Value = 3546
$ws.Range("I72").Value = 3244.3333
$ws.Range("K72").Value = 29198.9997
$ws.Range("M72").Value = -24830.9997

$ws.Range("H112").Value = 2639.0417
$ws.Range("J112").Value = 2639.0417
$ws.Range("L112").Value = 7917.125100000001
$ws.Range("N112").Value = -10133.1251

$ws.Range("H127").Value = 1124.5
$ws.Range("I127").Value = 853.5
$ws.Range("J127").Value = 1214.8334
$ws.Range("K127").Value = 2560.5
$ws.Range("L127").Value = 3644.5002
$ws.Range("M127").Value = 2399.5
$ws.Range("N127").Value = -13564.5002

$ws.Range("H129").Value = 883.6042
$ws.Range("J129").Value = 919.5
$ws.Range("L129").Value = 2758.5
$ws.Range("N129").Value = -12758.5

$ws.Range("H138").Value = 1260.5253
$ws.Range("I138").Value = 655.3333
$ws.Range("J138").Value = 1653.9
$ws.Range("K138").Value = 1965.9999
$ws.Range("L138").Value = 4961.700000000001
$ws.Range("M138").Value = 3174.0001
$ws.Range("N138").Value = -15241.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 261.25
$ws.Range("I4").Value = 423
$ws.Range("J4").Value = 99.5
$ws.Range("K4").Value = 423
$ws.Range("L4").Value = 99.5
$ws.Range("M4").Value = -307
$ws.Range("N4").Value = -331.5

$ws.Range("H5").Value = 221.57143
$ws.Range("I5").Value = 150.2
$ws.Range("K5").Value = 150.2
$ws.Range("M5").Value = -38.19999999999999

$ws.Range("H32").Value = 3544.4167
$ws.Range("I32").Value = 3655.121
$ws.Range("K32").Value = 3655.121
$ws.Range("M32").Value = -3368.121

$ws.Range("H74").Value = 1790.7
$ws.Range("I74").Value = 984.875
$ws.Range("J74").Value = 5014
$ws.Range("K74").Value = 984.875
$ws.Range("L74").Value = 5014
$ws.Range("M74").Value = -110.875
$ws.Range("N74").Value = -6762

$ws.Range("H77").Value = 1790.7
$ws.Range("I77").Value = 984.875
$ws.Range("J77").Value = 5014
$ws.Range("K77").Value = 4924.375
$ws.Range("L77").Value = 25070
$ws.Range("M77").Value = -556.375
$ws.Range("N77").Value = -33806

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H122").Value = 2142.8572
$ws.Range("I122").Value = 2083.3333
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 6249.999899999999
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -3799.999899999999
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 221.57143
$ws.Range("I4").Value = 150.2
$ws.Range("K4").Value = 150.2
$ws.Range("M4").Value = -35.19999999999999

$ws.Range("H64").Value = 320.83334
$ws.Range("I64").Value = 320
$ws.Range("J64").Value = 321.66666
$ws.Range("K64").Value = 320
$ws.Range("L64").Value = 321.66666
$ws.Range("M64").Value = -95
$ws.Range("N64").Value = -771.66666

$ws.Range("H67").Value = 320.83334
$ws.Range("I67").Value = 320
$ws.Range("J67").Value = 321.66666
$ws.Range("K67").Value = 320
$ws.Range("L67").Value = 321.66666
$ws.Range("M67").Value = 460
$ws.Range("N67").Value = -1881.66666

$ws.Range("H86").Value = 3971.6
$ws.Range("I86").Value = 4519.643
$ws.Range("J86").Value = 2692.8333
$ws.Range("K86").Value = 4519.643
$ws.Range("L86").Value = 2692.8333
$ws.Range("M86").Value = -3396.643
$ws.Range("N86").Value = -4938.8333

$ws.Range("H89").Value = 3971.6
$ws.Range("I89").Value = 4519.643
$ws.Range("J89").Value = 2692.8333
$ws.Range("K89").Value = 22598.215
$ws.Range("L89").Value = 13464.1665
$ws.Range("M89").Value = -16982.215
$ws.Range("N89").Value = -24696.1665

$ws.Range("H105").Value = 58825796
$ws.Range("I105").Value = 71430776
$ws.Range("J105").Value = 2543
$ws.Range("K105").Value = 71430776
$ws.Range("L105").Value = 2543
$ws.Range("M105").Value = -71429029
$ws.Range("N105").Value = -6037

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 222.90909
$ws.Range("I7").Value = 298.85715
$ws.Range("J7").Value = 90
$ws.Range("K7").Value = 298.85715
$ws.Range("L7").Value = 90
$ws.Range("M7").Value = -185.85715
$ws.Range("N7").Value = -316

$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H31").Value = 2134.6128
$ws.Range("I31").Value = 1105.5714
$ws.Range("J31").Value = 2982.0588
$ws.Range("K31").Value = 1105.5714
$ws.Range("L31").Value = 2982.0588
$ws.Range("M31").Value = -810.5714
$ws.Range("N31").Value = -3572.0588

$ws.Range("H34").Value = 2134.6128
$ws.Range("I34").Value = 1105.5714
$ws.Range("J34").Value = 2982.0588
$ws.Range("K34").Value = 1105.5714
$ws.Range("L34").Value = 2982.0588
$ws.Range("M34").Value = -903.5714
$ws.Range("N34").Value = -3386.0588

$ws.Range("H58").Value = 1357
$ws.Range("I58").Value = 1383.875
$ws.Range("J58").Value = 1249.5
$ws.Range("K58").Value = 1383.875
$ws.Range("L58").Value = 1249.5
$ws.Range("M58").Value = -1180.875
$ws.Range("N58").Value = -1655.5

$ws.Range("H99").Value = 1947.1428
$ws.Range("I99").Value = 1935
$ws.Range("K99").Value = 1935
$ws.Range("M99").Value = -437

$ws.Range("H122").Value = 950
$ws.Range("I122").Value = 840
$ws.Range("K122").Value = 2520
$ws.Range("M122").Value = -70

$ws.Range("H126").Value = 1947.1428
$ws.Range("I126").Value = 1935
$ws.Range("K126").Value = 5805
$ws.Range("M126").Value = -3335

$ws.Range("H136").Value = 1357
$ws.Range("I136").Value = 1383.875
$ws.Range("J136").Value = 1249.5
$ws.Range("K136").Value = 4151.625
$ws.Range("L136").Value = 3748.5
$ws.Range("M136").Value = -1601.625
$ws.Range("N136").Value = -8848.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 5550.1816
$ws.Range("I104").Value = 4684
$ws.Range("J104").Value = 5875
$ws.Range("K104").Value = 14052
$ws.Range("L104").Value = 17625
$ws.Range("M104").Value = -11431
$ws.Range("N104").Value = -22867

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

$ws.Range("H131").Value = 47620908
$ws.Range("I131").Value = 142857360
$ws.Range("J131").Value = 2680.1428
$ws.Range("K131").Value = 428572080
$ws.Range("L131").Value = 8040.428400000001
$ws.Range("M131").Value = -428567040
$ws.Range("N131").Value = -18120.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 4850
$ws.Range("J23").Value = 4850
$ws.Range("L23").Value = 4850
$ws.Range("N23").Value = -5296

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 232905.34
$ws.Range("J2").Value = 80444.53999999999
$ws.Range("L2").Value = 80444.53999999999
$ws.Range("N2").Value = -80668.53999999999

$ws.Range("H22").Value = 1198.85
$ws.Range("I22").Value = 1181.25
$ws.Range("J22").Value = 1225.25
$ws.Range("K22").Value = 1181.25
$ws.Range("L22").Value = 1225.25
$ws.Range("M22").Value = -886.25
$ws.Range("N22").Value = -1815.25

$ws.Range("H27").Value = 1198.85
$ws.Range("I27").Value = 1181.25
$ws.Range("J27").Value = 1225.25
$ws.Range("K27").Value = 1181.25
$ws.Range("L27").Value = 1225.25
$ws.Range("M27").Value = -1074.25
$ws.Range("N27").Value = -1439.25

$ws.Range("H46").Value = 2160.6
$ws.Range("I46").Value = 1767
$ws.Range("J46").Value = 2751
$ws.Range("K46").Value = 1767
$ws.Range("L46").Value = 2751
$ws.Range("M46").Value = -1579
$ws.Range("N46").Value = -3127

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3500
$ws.Range("I2").Value = 5000
$ws.Range("K2").Value = 5000
$ws.Range("M2").Value = -4888

$ws.Range("H122").Value = 54784936
$ws.Range("I122").Value = 63002532
$ws.Range("J122").Value = 966.6667
$ws.Range("K122").Value = 189007596
$ws.Range("L122").Value = 2900.0001
$ws.Range("M122").Value = -189005146
$ws.Range("N122").Value = -7800.0001
